$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Z34 cell (B4): append "Juju Worker OK" to the nickname line
$ws.Range("B4").Value = 'Z34: "Cosmic Cougar" Juju Worker OK' + "`n" + 'Worker _2_Z34 --> Ready' + "`n" + 'Host:192.168.1.61'

# Z10 cell (E4): split the comma-joined worker list into one "Juju Worker OK" line per worker
$ws.Range("E4").Value = 'Z10: ' + "`n" + 'Controladores K8s --> Ready' + "`n" + '"merry-tetra" Juju Worker OK ' + "`n" + '"pro-moose" Juju Worker OK' + "`n" + '"pumped-thrush" Juju Worker OK' + "`n" + 'Controlador Juju --> Ready' + "`n" + 'Host: 192.168.1.50'

# Z42 cell (C6): append "Juju Worker OK" and prefix the IP line with "Host "
$ws.Range("C6").Value = 'Z42:  "brave-wolf" Juju Worker OK' + "`n" + 'Worker_3_Z42--> Ready' + "`n" + 'Host 192.168.1.62'

# Z22 cell (E6): append "Juju Worker OK" to the nickname line
$ws.Range("E6").Value = 'Z22: "Exact-cicada" Juju Worker OK' + "`n" + 'Worker_8_Z22 --> Ready' + "`n" + 'Host: 192.168.1.57'

# Row 4 grew taller once the K8s worker cell wraps to more lines
$ws.Rows.Item(4).RowHeight = 118.5

# Selection moved from C7 to C6
$ws.Range("C6").Select() | Out-Null
